# Regenerate the "K" column (column G) values for the save_data sheet.
# The commit replaces the old "Strike#" derived values with recalculated
# K values for each row of match data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values keyed by data row number (row 2 .. row 13 on the sheet)
$kValues = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 1
    6  = 0
    7  = 0
    8  = 2
    9  = 3
    10 = 0
    11 = 1
    12 = 3
    13 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
